$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Feb 24 22:16:30 EST 2025"
$ws.Range("B3").Value = "Mon Feb 24 22:17:43 EST 2025"
$ws.Range("B4").Value = "Mon Feb 24 22:18:56 EST 2025"
$ws.Range("B5").Value = "Mon Feb 24 22:20:08 EST 2025"
$ws.Range("B6").Value = "Mon Feb 24 22:20:18 EST 2025"
$ws.Range("B7").Value = "Mon Feb 24 22:20:30 EST 2025"
$ws.Range("B8").Value = "Mon Feb 24 22:21:43 EST 2025"
$ws.Range("B9").Value = "Mon Feb 24 22:22:56 EST 2025"
$ws.Range("B10").Value = "Mon Feb 24 22:24:08 EST 2025"
$ws.Range("B11").Value = "Mon Feb 24 22:25:21 EST 2025"
$ws.Range("B12").Value = "Mon Feb 24 22:25:32 EST 2025"
$ws.Range("B13").Value = "Mon Feb 24 22:25:42 EST 2025"
$ws.Range("B14").Value = "Mon Feb 24 22:25:53 EST 2025"
$ws.Range("B15").Value = "Mon Feb 24 22:26:03 EST 2025"
$ws.Range("B16").Value = "Mon Feb 24 22:26:14 EST 2025"
$ws.Range("B17").Value = "Mon Feb 24 22:26:26 EST 2025"
$ws.Range("B18").Value = "Mon Feb 24 22:26:37 EST 2025"
$ws.Range("B19").Value = "Mon Feb 24 22:26:48 EST 2025"
$ws.Range("B20").Value = "Mon Feb 24 22:27:04 EST 2025"
$ws.Range("B21").Value = "Mon Feb 24 22:27:17 EST 2025"
$ws.Range("B22").Value = "Mon Feb 24 22:27:28 EST 2025"
$ws.Range("B23").Value = "Mon Feb 24 22:27:39 EST 2025"
$ws.Range("B24").Value = "Mon Feb 24 22:27:49 EST 2025"
$ws.Range("B25").Value = "Mon Feb 24 22:28:00 EST 2025"
$ws.Range("B26").Value = "Mon Feb 24 22:28:11 EST 2025"
$ws.Range("B27").Value = "Mon Feb 24 22:28:22 EST 2025"
$ws.Range("B28").Value = "Mon Feb 24 22:28:32 EST 2025"
$ws.Range("B29").Value = "Mon Feb 24 22:28:42 EST 2025"
$ws.Range("B30").Value = "Mon Feb 24 22:28:53 EST 2025"
$ws.Range("B31").Value = "Mon Feb 24 22:30:05 EST 2025"
$ws.Range("B32").Value = "Mon Feb 24 22:31:18 EST 2025"
$ws.Range("B33").Value = "Mon Feb 24 22:32:31 EST 2025"
$ws.Range("B34").Value = "Mon Feb 24 22:32:42 EST 2025"
$ws.Range("B35").Value = "Mon Feb 24 22:33:55 EST 2025"
$ws.Range("B36").Value = "Mon Feb 24 22:34:06 EST 2025"
$ws.Range("B37").Value = "Mon Feb 24 22:34:16 EST 2025"
$ws.Range("B38").Value = "Mon Feb 24 22:34:27 EST 2025"
$ws.Range("B39").Value = "Mon Feb 24 22:34:39 EST 2025"
$ws.Range("B40").Value = "Mon Feb 24 22:34:50 EST 2025"
$ws.Range("B41").Value = "Mon Feb 24 22:35:01 EST 2025"
$ws.Range("B42").Value = "Mon Feb 24 22:35:12 EST 2025"
$ws.Range("B43").Value = "Mon Feb 24 22:35:22 EST 2025"
$ws.Range("B44").Value = "Mon Feb 24 22:35:32 EST 2025"
$ws.Range("B45").Value = "Mon Feb 24 22:35:43 EST 2025"
$ws.Range("B46").Value = "Mon Feb 24 22:35:54 EST 2025"
$ws.Range("B47").Value = "Mon Feb 24 22:36:06 EST 2025"
$ws.Range("B48").Value = "Mon Feb 24 22:36:17 EST 2025"
$ws.Range("B49").Value = "Mon Feb 24 22:36:28 EST 2025"
$ws.Range("B50").Value = "Mon Feb 24 22:37:41 EST 2025"
$ws.Range("B51").Value = "Mon Feb 24 22:37:52 EST 2025"
$ws.Range("B52").Value = "Mon Feb 24 22:38:03 EST 2025"
$ws.Range("B53").Value = "Mon Feb 24 22:38:14 EST 2025"
